# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet
#    and populate it with the quarterly fund-holding detail rows.
# 2. Insert a new leading data row into "总计" summarizing the 2022-Q1
#    quarter (count of holdings + total held value), shifting the existing
#    rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: apply the bold / thin-box-border / center-top style that is
# used for header cells and the row-index column throughout this
# workbook.
# ---------------------------------------------------------------------
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $cell.Borders.Item(9).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $cell.Borders.Item(10).Weight = 2
}

# Helper: write a value into a cell forcing it to be stored as literal
# text (so numeric-looking strings like "5.02" are not reinterpreted as
# numbers / floats).
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet just before "总计"
# ---------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheetBefore)
$q1.Name = "2022-Q1"

# Header row (row 1, columns B..H)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q1.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    Set-HeaderStyle($cell)
}

# Data rows (rows 2..9): A=index, B=code, C=name, D=scale, E=stock pos,
# F=pos pct, G=held value, H=rank
$rows = @(
    @("002181", "华安大安全主题灵活配置混合",             "5.02", "87.28", "3.03", "0.1521", 9),
    @("012084", "博时睿弘一年定期开放混合型证券投资基金A", "1.78", "91.29", "4.26", "0.0758", 3),
    @("519113", "浦银安盛精致生活混合",                   "2.09", "90.20", "2.45", "0.0512", 9),
    @("010194", "博时睿祥15个月定期开放混合A",             "0.42", "79.12", "6.97", "0.0293", 4),
    @("001273", "民生加银新动力灵活配置混合A",             "0.04", "68.44", "3.63", "0.0015", 1),
    @("001274", "民生加银新动力灵活配置混合D",             "0.04", "68.44", "3.63", "0.0015", 1),
    @("010195", "博时睿祥15个月定期开放混合C",             "0.00", "79.12", "6.97", $null,     4),
    @("012085", "博时睿弘一年定期开放混合型证券投资基金C", "0.00", "91.29", "4.26", $null,     3)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $idxCell = $q1.Cells.Item($r, 1)
    $idxCell.Value = $i
    Set-HeaderStyle($idxCell)

    Set-TextValue $q1.Cells.Item($r, 2) $data[0]
    Set-TextValue $q1.Cells.Item($r, 3) $data[1]
    Set-TextValue $q1.Cells.Item($r, 4) $data[2]
    Set-TextValue $q1.Cells.Item($r, 5) $data[3]
    Set-TextValue $q1.Cells.Item($r, 6) $data[4]

    if ($null -eq $data[5]) {
        $q1.Cells.Item($r, 7).Value = 0
    } else {
        Set-TextValue $q1.Cells.Item($r, 7) $data[5]
    }

    $q1.Cells.Item($r, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# Step 2: insert the new 2022-Q1 summary row at the top of "总计"
#
# NOTE: the worksheet collection re-indexes when a new sheet is
# inserted, so we must look the "总计" sheet up again by name instead
# of reusing the handle obtained before the insert.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$a2 = $totalSheet.Cells.Item(2, 1)
$a2.Value = 0
Set-HeaderStyle($a2)

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 0.31

# Renumber the index column for the rows that shifted down
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
